$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

function Clear-Cell($ws, $cellRef) {
    $ws.Range($cellRef).ClearContents()
}


# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17
Set-Cell $ws "H17" 1014.53845
Set-Cell $ws "I17" 540
Set-Cell $ws "J17" 1311.125
Set-Cell $ws "K17" 1620
Set-Cell $ws "L17" 3933.375
Set-Cell $ws "M17" -1452
Set-Cell $ws "N17" -4269.375

# Row 51
Set-Cell $ws "H51" 9094999
Set-Cell $ws "I51" 18186618
Set-Cell $ws "J51" 3380
Set-Cell $ws "K51" 18186618
Set-Cell $ws "L51" 3380
Set-Cell $ws "M51" -18186134
Set-Cell $ws "N51" -4348

# Row 70
Set-Cell $ws "H70" 1370.44
Set-Cell $ws "I70" 1359.174
Set-Cell $ws "J70" 1500
Set-Cell $ws "K70" 4077.522
Set-Cell $ws "L70" 4500
Set-Cell $ws "M70" -3807.522
Set-Cell $ws "N70" -5040

# Row 73
Set-Cell $ws "H73" 1370.44
Set-Cell $ws "I73" 1359.174
Set-Cell $ws "J73" 1500
Set-Cell $ws "K73" 4077.522
Set-Cell $ws "L73" 4500
Set-Cell $ws "M73" -3141.522
Set-Cell $ws "N73" -6372

# Row 112
Set-Cell $ws "H112" 3408.3333
Set-Cell $ws "I112" 1500
Set-Cell $ws "J112" 3910.5264
Set-Cell $ws "K112" 4500
Set-Cell $ws "L112" 11731.5792
Set-Cell $ws "M112" -3392
Set-Cell $ws "N112" -13947.5792

# Row 113
Set-Cell $ws "H113" 2961.7
Set-Cell $ws "I113" 2442.8333
Set-Cell $ws "J113" 3740
Set-Cell $ws "K113" 2442.8333
Set-Cell $ws "L113" 3740
Set-Cell $ws "M113" 811.1667000000002
Set-Cell $ws "N113" -10248

# Row 127
Set-Cell $ws "H127" 1315
Set-Cell $ws "I127" 894.5
Set-Cell $ws "J127" 1391.4546
Set-Cell $ws "K127" 2683.5
Set-Cell $ws "L127" 4174.3638
Set-Cell $ws "M127" 2276.5
Set-Cell $ws "N127" -14094.3638

# Row 129
Set-Cell $ws "H129" 2045.36
Set-Cell $ws "I129" 812.4
Set-Cell $ws "J129" 2353.6
Set-Cell $ws "K129" 2437.2
Set-Cell $ws "L129" 7060.799999999999
Set-Cell $ws "M129" 2562.8
Set-Cell $ws "N129" -17060.8

# Row 138
Set-Cell $ws "H138" 2937.0896
Set-Cell $ws "I138" 1935.6666
Set-Cell $ws "J138" 3613.05
Set-Cell $ws "K138" 5806.9998
Set-Cell $ws "L138" 10839.15
Set-Cell $ws "M138" -666.9997999999996
Set-Cell $ws "N138" -21119.15


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 24
Set-Cell $ws "H24" 0
Set-Cell $ws "I24" 0
Set-Cell $ws "J24" 0
Set-Cell $ws "K24" 0
Set-Cell $ws "L24" 0
Clear-Cell $ws "N24"

# Row 32
Set-Cell $ws "H32" 5942.3
Set-Cell $ws "I32" 5705.521
Set-Cell $ws "J32" 11625
Set-Cell $ws "K32" 5705.521
Set-Cell $ws "L32" 11625
Set-Cell $ws "M32" -5418.521
Set-Cell $ws "N32" -12199

# Row 76
Set-Cell $ws "H76" 14766.667
Set-Cell $ws "I76" 0
Set-Cell $ws "J76" 14766.667
Set-Cell $ws "K76" 0
Set-Cell $ws "L76" 14766.667
Set-Cell $ws "N76" -15442.667

# Row 79
Set-Cell $ws "H79" 14766.667
Set-Cell $ws "I79" 0
Set-Cell $ws "J79" 14766.667
Set-Cell $ws "K79" 0
Set-Cell $ws "L79" 14766.667
Set-Cell $ws "N79" -17106.667

# Row 80
Set-Cell $ws "H80" 40000
Set-Cell $ws "I80" 0
Set-Cell $ws "J80" 40000
Set-Cell $ws "K80" 0
Set-Cell $ws "L80" 40000
Set-Cell $ws "N80" -41996

# Row 83
Set-Cell $ws "H83" 40000
Set-Cell $ws "I83" 0
Set-Cell $ws "J83" 40000
Set-Cell $ws "K83" 0
Set-Cell $ws "L83" 120000
Set-Cell $ws "N83" -129984

# Row 100
Set-Cell $ws "H100" 0
Set-Cell $ws "I100" 0
Set-Cell $ws "J100" 0
Set-Cell $ws "K100" 0
Set-Cell $ws "L100" 0
Clear-Cell $ws "N100"

# Row 122
Set-Cell $ws "H122" 1373.1
Set-Cell $ws "I122" 1157.6842
Set-Cell $ws "J122" 1745.1818
Set-Cell $ws "K122" 3473.0526
Set-Cell $ws "L122" 5235.5454
Set-Cell $ws "M122" -1023.0526
Set-Cell $ws "N122" -10135.5454

# Row 132
Set-Cell $ws "H132" 6279.6
Set-Cell $ws "I132" 6391.6924
Set-Cell $ws "J132" 6158.1665
Set-Cell $ws "K132" 19175.0772
Set-Cell $ws "L132" 18474.4995
Set-Cell $ws "M132" -16645.0772
Set-Cell $ws "N132" -23534.4995


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 54
Set-Cell $ws "H54" 0
Set-Cell $ws "I54" 0
Set-Cell $ws "J54" 0
Set-Cell $ws "K54" 0
Set-Cell $ws "L54" 0
Clear-Cell $ws "M54"
Clear-Cell $ws "N54"


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 3
Set-Cell $ws "H3" 2999.5
Set-Cell $ws "I3" 0
Set-Cell $ws "J3" 2999.5
Set-Cell $ws "K3" 0
Set-Cell $ws "L3" 2999.5
Set-Cell $ws "N3" -3225.5
Clear-Cell $ws "M3"

# Row 5
Set-Cell $ws "H5" 376
Set-Cell $ws "I5" 150.125
Set-Cell $ws "J5" 737.4
Set-Cell $ws "K5" 150.125
Set-Cell $ws "L5" 737.4
Set-Cell $ws "M5" -38.125
Set-Cell $ws "N5" -961.4

# Row 22
Set-Cell $ws "H22" 229.125
Set-Cell $ws "I22" 233.28572
Set-Cell $ws "J22" 200
Set-Cell $ws "K22" 233.28572
Set-Cell $ws "L22" 200
Set-Cell $ws "M22" 116.71428
Set-Cell $ws "N22" -900

# Row 62
Set-Cell $ws "H62" 14718.75
Set-Cell $ws "I62" 2535.7144
Set-Cell $ws "J62" 100000
Set-Cell $ws "K62" 2535.7144
Set-Cell $ws "L62" 100000
Set-Cell $ws "M62" -1911.7144
Set-Cell $ws "N62" -101248

# Row 65
Set-Cell $ws "H65" 14718.75
Set-Cell $ws "I65" 2535.7144
Set-Cell $ws "J65" 100000
Set-Cell $ws "K65" 12678.572
Set-Cell $ws "L65" 500000
Set-Cell $ws "M65" -9558.572
Set-Cell $ws "N65" -506240

# Row 86
Set-Cell $ws "H86" 4773.154
Set-Cell $ws "I86" 4116.4736
Set-Cell $ws "J86" 6555.5713
Set-Cell $ws "K86" 4116.4736
Set-Cell $ws "L86" 6555.5713
Set-Cell $ws "M86" -2993.4736
Set-Cell $ws "N86" -8801.5713

# Row 89
Set-Cell $ws "H89" 4773.154
Set-Cell $ws "I89" 4116.4736
Set-Cell $ws "J89" 6555.5713
Set-Cell $ws "K89" 20582.368
Set-Cell $ws "L89" 32777.85649999999
Set-Cell $ws "M89" -14966.368
Set-Cell $ws "N89" -44009.85649999999

# Row 92
Set-Cell $ws "H92" 14368.2
Set-Cell $ws "I92" 0
Set-Cell $ws "J92" 14368.2
Set-Cell $ws "K92" 0
Set-Cell $ws "L92" 14368.2
Set-Cell $ws "N92" -19360.2


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 3
Set-Cell $ws "H3" 6739.8667
Set-Cell $ws "I3" 1891.25
Set-Cell $ws "J3" 12281.143
Set-Cell $ws "K3" 5673.75
Set-Cell $ws "L3" 36843.429
Set-Cell $ws "M3" -5561.75
Set-Cell $ws "N3" -37067.429

# Row 9
Set-Cell $ws "H9" 335666.66
Set-Cell $ws "I9" 500750
Set-Cell $ws "J9" 5500
Set-Cell $ws "K9" 1502250
Set-Cell $ws "L9" 16500
Set-Cell $ws "M9" -1502026
Set-Cell $ws "N9" -16948

# Row 113
Set-Cell $ws "H113" 784.7059
Set-Cell $ws "I113" 699
Set-Cell $ws "J113" 881.125
Set-Cell $ws "K113" 2097
Set-Cell $ws "L113" 2643.375
Set-Cell $ws "M113" 73
Set-Cell $ws "N113" -6983.375

# Row 131
Set-Cell $ws "H131" 1412.4445
Set-Cell $ws "I131" 1395
Set-Cell $ws "J131" 1419.1538
Set-Cell $ws "K131" 4185
Set-Cell $ws "L131" 4257.4614
Set-Cell $ws "M131" 855
Set-Cell $ws "N131" -14337.4614

# Row 133
Set-Cell $ws "H133" 4403.3335
Set-Cell $ws "I133" 1921.4286
Set-Cell $ws "J133" 5071.5386
Set-Cell $ws "K133" 5764.2858
Set-Cell $ws "L133" 15214.6158
Set-Cell $ws "M133" -704.2857999999997
Set-Cell $ws "N133" -25334.6158


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
Set-Cell $ws "H122" 1856.6552
Set-Cell $ws "I122" 1732.1765
Set-Cell $ws "J122" 2033
Set-Cell $ws "K122" 5196.529500000001
Set-Cell $ws "L122" 6099
Set-Cell $ws "M122" -2746.529500000001
Set-Cell $ws "N122" -10999


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
Set-Cell $ws "H2" 2311997
Set-Cell $ws "I2" 1000
Set-Cell $ws "J2" 2504580
Set-Cell $ws "K2" 1000
Set-Cell $ws "L2" 2504580
Set-Cell $ws "M2" -888
Set-Cell $ws "N2" -2504804

# Row 139
Set-Cell $ws "H139" 45000
Set-Cell $ws "I139" 30000
Set-Cell $ws "J139" 50000
Set-Cell $ws "K139" 30000
Set-Cell $ws "L139" 50000
Set-Cell $ws "M139" -24860
Set-Cell $ws "N139" -60280


# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 136
Set-Cell $ws "H136" 1845.6182
Set-Cell $ws "I136" 1856.5745
Set-Cell $ws "J136" 1781.25
Set-Cell $ws "K136" 5569.7235
Set-Cell $ws "L136" 5343.75
Set-Cell $ws "M136" -3019.7235
Set-Cell $ws "N136" -10443.75
